# Updates cryptos price list cells (Price/Volume columns) to match the
# latest scrape, mirroring the GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.237.14"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.588.50"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'211.91"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").Value = "'0.503"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.244"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "'19.21"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "1.812.69"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "1.582.00"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "'0.514"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").Value = "'63.86"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "26.245.15"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "0.0₃0725"
$ws.Range("D19").Value = "'7.45"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "'214.16"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").Value = "'8.98"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").Value = "'144.24"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").Value = "'15.09"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("E30").Value = "  -2.02%  "
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("D33").Value = "1.416.19"
$ws.Range("E33").Value = "  +7.97%  "
$ws.Range("D34").Value = "'2.94"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("D36").Value = "'0.590"
$ws.Range("E36").Value = "  -4.30%  "
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("E39").Value = "  +5.12%  "
$ws.Range("D40").Value = "'0.822"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'0.940"
$ws.Range("E42").Value = "  -14.45%  "
$ws.Range("D43").Value = "'0.766"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "1.723.90"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("D46").Value = "'61.14"
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("D47").Value = "'85.78"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'0.0502"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "'0.0967"
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.03%  "
